$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 81; $r++) {
    $ws.Cells.Item($r, 3).Value = 0
    $ws.Cells.Item($r, 4).Value = 0
    $ws.Cells.Item($r, 6).Value = 0
}

$ws.Columns.Item(6).AutoFit()
